# "Add files via upload" -- re-exported workbook with an updated CMS label
# and a different last-used cell selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The L column (L2:L7) holds a shared "CMS(<mu>)" label; append "(2010)" to
# it so it reads "CMS(<mu>)(2010)". Writing the same text to every cell that
# shares the string keeps them consolidated on a single shared-string entry.
$newLabel = "CMS(" + [char]0x00B5 + ")(2010)"
$ws.Range("L2:L7").Value = $newLabel

# Update the sheet's last-used selection.
[void]$ws.Range("L13").Select()
